# Updates the SeenRx_CKJ sheet:
#  - Strips the " Seen Rx" suffix from the FFTR header labels (and upper-cases NOCLOG)
#  - Replaces the CKJ product rows with the full, sequential CKJ product list
#    (CKJ, CKJ10-CKJ16, CKJ20-CKJ26, CKJ30-CKJ35, CKJ40-CKJ46, CKJ50-CKJ56)
#    and refreshes the counts for each column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1, columns B-J); column A (FFTR) is unchanged ----
$headers = @("LIGAZID","EMAZID","LIPICON","AGLIP","CIFIBET","AMLEVO","CARDOBIS","RIVAROX","NOCLOG")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---- Data rows: row number, label, then the 9 counts for columns B..J ----
$data = @(
    ,@(2, "CKJ", 44,19,22,7,10,22,41,5,54)
    ,@(3, "CKJ10", 14,3,6,1,4,6,14,1,10)
    ,@(4, "CKJ11", 1,1,2,0,0,1,3,0,4)
    ,@(5, "CKJ12", 7,1,1,0,2,1,3,0,0)
    ,@(6, "CKJ13", 2,0,3,1,1,1,3,0,3)
    ,@(7, "CKJ14", 3,0,0,0,0,0,1,1,1)
    ,@(8, "CKJ15", 1,0,0,0,1,2,2,0,2)
    ,@(9, "CKJ16", 0,1,0,0,0,1,2,0,0)
    ,@(10, "CKJ20", 9,6,5,0,2,1,10,0,15)
    ,@(11, "CKJ21", 0,0,0,0,0,0,2,0,0)
    ,@(12, "CKJ22", 4,5,0,0,0,1,1,0,0)
    ,@(13, "CKJ23", 0,0,2,0,0,0,4,0,11)
    ,@(14, "CKJ24", 5,1,2,0,2,0,3,0,3)
    ,@(15, "CKJ25", 0,0,1,0,0,0,0,0,1)
    ,@(16, "CKJ26", 0,0,0,0,0,0,0,0,0)
    ,@(17, "CKJ30", 11,7,7,4,4,7,9,3,17)
    ,@(18, "CKJ31", 3,4,0,1,1,3,2,0,3)
    ,@(19, "CKJ32", 2,1,7,1,0,1,5,1,12)
    ,@(20, "CKJ33", 6,1,0,2,0,2,2,2,0)
    ,@(21, "CKJ34", 0,0,0,0,0,0,0,0,0)
    ,@(22, "CKJ35", 0,1,0,0,3,1,0,0,2)
    ,@(23, "CKJ40", 4,2,1,0,0,6,4,1,11)
    ,@(24, "CKJ41", 0,0,0,0,0,0,0,0,0)
    ,@(25, "CKJ42", 0,0,0,0,0,3,2,1,10)
    ,@(26, "CKJ43", 0,0,0,0,0,0,0,0,0)
    ,@(27, "CKJ44", 3,2,0,0,0,1,0,0,0)
    ,@(28, "CKJ45", 0,0,0,0,0,0,0,0,0)
    ,@(29, "CKJ46", 1,0,1,0,0,2,2,0,1)
    ,@(30, "CKJ50", 6,1,3,2,0,2,4,0,1)
    ,@(31, "CKJ51", 1,0,0,1,0,1,1,0,0)
    ,@(32, "CKJ52", 0,0,0,0,0,0,0,0,0)
    ,@(33, "CKJ53", 0,0,0,0,0,0,0,0,0)
    ,@(34, "CKJ54", 2,0,1,1,0,0,2,0,1)
    ,@(35, "CKJ55", 3,1,2,0,0,1,1,0,0)
    ,@(36, "CKJ56", 0,0,0,0,0,0,0,0,0)
)

foreach ($row in $data) {
    $r = $row[0]
    $label = $row[1]
    $ws.Cells.Item($r, 1).Value = $label
    for ($c = 0; $c -lt 9; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $row[2 + $c]
    }
}

Write-Host ("Done. UsedRange=" + $ws.UsedRange.Address())
